# Adds a new "2022-Q3" quarterly sheet (with its fund-holding detail data)
# ahead of the existing "2022-Q2" sheet, and updates the "总计" (summary)
# sheet with a new top row for 2022-Q3, pushing the older quarters down.

$wb = $excel.ActiveWorkbook

# Helper: force a value to be stored as literal text (so numeric-looking
# strings like "004823" or "5.35" keep their leading zeros / exact text
# instead of being auto-coerced into numbers), then drop back to the
# default "Normal" style so no stray number-format style sticks around.
function Set-TextCell {
    param($sheet, $addr, $val)
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" worksheet right after "总计" (so the
# tab order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4).
#
# NOTE: worksheet references must be (re-)fetched *after* Add() runs -
# handles grabbed beforehand can end up pointing at stale sheet data
# once the collection is mutated, which silently drops pasted styles.
# ---------------------------------------------------------------------
$totalSheetTmp = $wb.Worksheets.Item("总计")
$null = $wb.Worksheets.Add($null, $totalSheetTmp)

$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"
$dstSheet = $wb.Worksheets.Item("2022-Q3")
$dstSheet.Activate()
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# Step 2: clone the formatting (header row + index column) from the
# "2022-Q2" sheet so the new sheet matches the existing house style.
# ---------------------------------------------------------------------
$q2Sheet.Range("B1:H1").Copy()
$dstSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q2Sheet.Range("A2").Copy()
$dstSheet.Range("A2:A7").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# Step 3: header row labels (plain text, no numeric look, keeps the
# bold/bordered style copied above).
# ---------------------------------------------------------------------
$dstSheet.Range("B1").Value = "基金代码"
$dstSheet.Range("C1").Value = "基金名称"
$dstSheet.Range("D1").Value = "基金规模"
$dstSheet.Range("E1").Value = "股票总仓位"
$dstSheet.Range("F1").Value = "仓位占比"
$dstSheet.Range("G1").Value = "持有市值(亿元)"
$dstSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# Step 4: fund-holding detail rows for 2022-Q3.
# Column A = row index (number), column H = rank (number);
# columns B-G are stored as literal text, matching the source data.
# ---------------------------------------------------------------------
$dstSheet.Range("A2").Value = 0
$dstSheet.Range("H2").Value = 7
Set-TextCell $dstSheet "B2" "004823"
Set-TextCell $dstSheet "C2" "上投摩根安裕回报混合A"
Set-TextCell $dstSheet "D2" "5.35"
Set-TextCell $dstSheet "E2" "25.71"
Set-TextCell $dstSheet "F2" "1.38"
Set-TextCell $dstSheet "G2" "0.0738"

$dstSheet.Range("A3").Value = 1
$dstSheet.Range("H3").Value = 7
Set-TextCell $dstSheet "B3" "004824"
Set-TextCell $dstSheet "C3" "上投摩根安裕回报混合C"
Set-TextCell $dstSheet "D3" "4.91"
Set-TextCell $dstSheet "E3" "25.71"
Set-TextCell $dstSheet "F3" "1.38"
Set-TextCell $dstSheet "G3" "0.0678"

$dstSheet.Range("A4").Value = 2
$dstSheet.Range("H4").Value = 7
Set-TextCell $dstSheet "B4" "001231"
Set-TextCell $dstSheet "C4" "银华泰利灵活配置混合A"
Set-TextCell $dstSheet "D4" "0.96"
Set-TextCell $dstSheet "E4" "24.09"
Set-TextCell $dstSheet "F4" "0.88"
Set-TextCell $dstSheet "G4" "0.0084"

$dstSheet.Range("A5").Value = 3
$dstSheet.Range("H5").Value = 5
Set-TextCell $dstSheet "B5" "003063"
Set-TextCell $dstSheet "C5" "银华通利灵活配置混合C"
Set-TextCell $dstSheet "D5" "0.30"
Set-TextCell $dstSheet "E5" "26.56"
Set-TextCell $dstSheet "F5" "1.19"
Set-TextCell $dstSheet "G5" "0.0036"

$dstSheet.Range("A6").Value = 4
$dstSheet.Range("H6").Value = 5
Set-TextCell $dstSheet "B6" "003062"
Set-TextCell $dstSheet "C6" "银华通利灵活配置混合A"
Set-TextCell $dstSheet "D6" "0.22"
Set-TextCell $dstSheet "E6" "26.56"
Set-TextCell $dstSheet "F6" "1.19"
Set-TextCell $dstSheet "G6" "0.0026"

$dstSheet.Range("A7").Value = 5
$dstSheet.Range("H7").Value = 7
Set-TextCell $dstSheet "B7" "002328"
Set-TextCell $dstSheet "C7" "银华泰利灵活配置混合C"
Set-TextCell $dstSheet "D7" "0.03"
Set-TextCell $dstSheet "E7" "24.09"
Set-TextCell $dstSheet "F7" "0.88"
Set-TextCell $dstSheet "G7" "0.0003"

# ---------------------------------------------------------------------
# Step 5: update the "总计" (summary) sheet - insert a new 2022-Q3 row
# at the top of the data (row 2) and push the other quarters down one
# row. Rewritten directly since every final value is already known.
# Row 5 is brand new, so clone A2's styling (index-column style) onto
# A5 before writing its value.
# ---------------------------------------------------------------------
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.16

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 0.02

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.03

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 6
$totalSheet.Range("D5").Value = 0.53
